$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (optimisation summary) ---
$wsSchedule.Range("E4").Value = 497.91622725
$wsSchedule.Range("F4").Value = 29.27197103174603
$wsSchedule.Range("E5").Value = 627.4465957500001
$wsSchedule.Range("F5").Value = 18.44346254409172

# --- Detailed sheet updates (price / type forecasts) ---
$wsDetailed.Range("B39").Value = 70.48375
$wsDetailed.Range("B40").Value = 101.47343
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 57.09
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 80.02
$wsDetailed.Range("B44").Value = 80.02
$wsDetailed.Range("B46").Value = 85.95
$wsDetailed.Range("B50").Value = 58.22268
$wsDetailed.Range("B51").Value = 58.41034
$wsDetailed.Range("B52").Value = 57.88255
$wsDetailed.Range("B53").Value = 56.98
$wsDetailed.Range("B54").Value = 51.26774
$wsDetailed.Range("B57").Value = 56.98
$wsDetailed.Range("B58").Value = 56.98
$wsDetailed.Range("B59").Value = 64.94638999999999
$wsDetailed.Range("B60").Value = 65
$wsDetailed.Range("B61").Value = 76.28136000000001
$wsDetailed.Range("B62").Value = 70.28973999999999
$wsDetailed.Range("B65").Value = 36.07
$wsDetailed.Range("B66").Value = 36.05989
$wsDetailed.Range("B70").Value = 41.72921
$wsDetailed.Range("B71").Value = 46.61275
$wsDetailed.Range("B72").Value = 36.06028
$wsDetailed.Range("B73").Value = 39.065
$wsDetailed.Range("B74").Value = 36.06
$wsDetailed.Range("B75").Value = 47.65777
$wsDetailed.Range("B77").Value = 28.81204
$wsDetailed.Range("B78").Value = 36.0601
$wsDetailed.Range("B79").Value = 36.05857
$wsDetailed.Range("B80").Value = 16.83806
$wsDetailed.Range("B81").Value = 26.1512
$wsDetailed.Range("B82").Value = 15.65567
$wsDetailed.Range("B83").Value = 6.48108
$wsDetailed.Range("B84").Value = 5.98882
$wsDetailed.Range("B85").Value = -9.555009999999999
$wsDetailed.Range("B86").Value = -6.88159
$wsDetailed.Range("B87").Value = -3.07461
$wsDetailed.Range("B88").Value = -3.07143
$wsDetailed.Range("B93").Value = 64.8901
$wsDetailed.Range("B94").Value = 57.09
